$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("branch")

# The previously blank row 6 only carried placeholder formatting
# (style "1" on A6). Reset the whole row to the plain "Normal" style
# used by the other data rows before filling it in.
$ws.Range("A6:D6").Style = "Normal"

# Add a new branch row (row 6) with data matching the other branch entries
# (id, name, location, staffQuota).
$ws.Range("A6").Value = "6b4b325a-a060-46ce-969c-a5a427566f4b"
$ws.Range("B6").Value = "NTU"
$ws.Range("C6").Value = "North Spine Plaza"

# staffQuota is stored as text (like row 4's "10"), so force text entry
# for D6, then restore the Normal style so the cell ends up unstyled
# just like the other data rows.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "10"
$ws.Range("D6").Style = "Normal"
